# Generate Report for Handoff
# Update the handoff id (guid) and related timestamps across the three sheets.

$wb = $excel.ActiveWorkbook

$oldId = "1a3aef7f-8dea-45ee-bab0-e6b816cef846"
$newId = "85579677-29ce-4be0-9469-998e9b3540b8"

$oldZhHash = "ec6af8b10ffa78461560440d303b57f0a1fd1cd7"
$newZhHash = "ec2286c49513a49fd9cda885914b69bda8b27d74"

$oldDeHash = "ec6af8b10ffa78461560440d303b57f0a1fd1cd7"
$newDeHash = "ec2286c49513a49fd9cda885914b69bda8b27d74"

$newMdName = "$newId.md"
$newZhName = "$newId.$newZhHash.zh-cn.xlf"
$newDeName = "$newId.$newDeHash.de-de.xlf"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("D2").Value = "2016-46-20 06:46:55"
foreach ($h in $wsOverview.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq "`$A`$2") {
        $h.TextToDisplay = $newMdName
    }
}

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("D2").Value = $newZhName
$wsZh.Range("E2").Value = "2016-03-20 06:46:52"
foreach ($h in $wsZh.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq "`$A`$2") {
        $h.TextToDisplay = $newMdName
    } elseif ($addr -eq "`$D`$2") {
        $h.TextToDisplay = $newZhName
    }
}

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("D2").Value = $newDeName
$wsDe.Range("E2").Value = "2016-03-20 06:46:55"
foreach ($h in $wsDe.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq "`$A`$2") {
        $h.TextToDisplay = $newMdName
    } elseif ($addr -eq "`$D`$2") {
        $h.TextToDisplay = $newDeName
    }
}
